$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: values from original row 5
$ws.Range("D2").Value = 44965
$ws.Range("J2").Value = 1120

# Row 3: values from original row 9
$ws.Range("D3").Value = 45084
$ws.Range("J3").Value = 900

# Row 4: values from original row 8
$ws.Range("D4").Value = 44911
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 1800
$ws.Range("M4").Value = 1900
$ws.Range("P4").Value = 633

# Row 5: values from original row 13
$ws.Range("D5").Value = 44985
$ws.Range("J5").Value = 1000

# Row 6: values from original row 30
$ws.Range("D6").Value = 45077
$ws.Range("J6").Value = 760

# Row 7: values from original row 4
$ws.Range("D7").Value = 44881
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1900
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1950
$ws.Range("P7").Value = 650

# Row 8: values from original row 18
$ws.Range("D8").Value = 44685
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 1500
$ws.Range("M8").Value = 1750
$ws.Range("P8").Value = 583

# Row 9: values from original row 7
$ws.Range("D9").Value = 45070
$ws.Range("J9").Value = 800

# Row 10: values from original row 21
$ws.Range("D10").Value = 44953
$ws.Range("J10").Value = 1000

# Row 11: values from original row 26
$ws.Range("D11").Value = 44992
$ws.Range("J11").Value = 1040

# Row 12: values from original row 31
$ws.Range("D12").Value = 45007
$ws.Range("J12").Value = 1160
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2250
$ws.Range("P12").Value = 750

# Row 13: values from original row 16
$ws.Range("D13").Value = 45034
$ws.Range("J13").Value = 1100

# Row 14: values from original row 19
$ws.Range("D14").Value = 45020
$ws.Range("J14").Value = 1200

# Row 15: values from original row 22
$ws.Range("D15").Value = 44951

# Row 16: values from original row 14
$ws.Range("D16").Value = 44970
$ws.Range("J16").Value = 800

# Row 17: values from original row 25
$ws.Range("D17").Value = 45062
$ws.Range("J17").Value = 1100
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2250
$ws.Range("P17").Value = 750

# Row 18: values from original row 27
$ws.Range("D18").Value = 44971
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2250
$ws.Range("P18").Value = 750

# Row 19: values from original row 2
$ws.Range("D19").Value = 45013
$ws.Range("J19").Value = 1100

# Row 20: values from original row 28
$ws.Range("D20").Value = 44848
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 1750
$ws.Range("P20").Value = 583

# Row 21: values from original row 15
$ws.Range("D21").Value = 45091
$ws.Range("J21").Value = 800

# Row 22: values from original row 20
$ws.Range("D22").Value = 45035
$ws.Range("J22").Value = 1100

# Row 24: values from original row 11
$ws.Range("D24").Value = 44827
$ws.Range("J24").Value = 1200
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2500
$ws.Range("M24").Value = 2250
$ws.Range("P24").Value = 750

# Row 25: values from original row 6
$ws.Range("D25").Value = 45006

# Row 26: values from original row 24
$ws.Range("D26").Value = 44978
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 1800
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 1900
$ws.Range("P26").Value = 633

# Row 27: values from original row 3
$ws.Range("D27").Value = 45028

# Row 28: values from original row 29
$ws.Range("D28").Value = 45041
$ws.Range("J28").Value = 1160
$ws.Range("K28").Value = 2000
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2250
$ws.Range("P28").Value = 750

# Row 29: values from original row 10
$ws.Range("D29").Value = 44999
$ws.Range("J29").Value = 1100

# Row 30: values from original row 12
$ws.Range("D30").Value = 44910
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 1800
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = 1900
$ws.Range("P30").Value = 633

# Row 31: values from original row 17
$ws.Range("D31").Value = 44883
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 1900
$ws.Range("P31").Value = 633
